# Apply the OF.xlsx edit:
#  1. Rename sheet "Chappe" -> "Chape"
#  2. Add a new sheet "Support de skate" at the end of the workbook
#  3. Fill in the new sheet's data
#  4. On the "Chape" sheet, change machine name "Four" -> "Fraiseuse 3 axes manuelle"
#     for rows 5-9, and fix row 9's numeric-looking entries to be stored as text
#     (matching the rest of the column, which is text)

$wb = $excel.ActiveWorkbook

# --- 1. Rename "Chappe" to "Chape" ---
$chape = $wb.Sheets.Item("Chappe")
$chape.Name = "Chape"

# --- 2. Update machine names on the Chape sheet (rows 5-8 stay as they were,
#         just swap the machine text) ---
$chape.Range("C5").Value = "Fraiseuse 3 axes manuelle"
$chape.Range("C6").Value = "Fraiseuse 3 axes manuelle"
$chape.Range("C7").Value = "Fraiseuse 3 axes manuelle"
$chape.Range("C8").Value = "Fraiseuse 3 axes manuelle"

# Row 9: keep values "5"/"1" but force them to remain text, like the rest
# of their columns (A5:A8 and D5:D8 are text already).
$chape.Range("A9").NumberFormat = "@"
$chape.Range("A9").Value = "5"

$chape.Range("C9").Value = "Fraiseuse 3 axes manuelle"

$chape.Range("D9").NumberFormat = "@"
$chape.Range("D9").Value = "1"

# --- 3. Add the new "Support de skate" sheet at the end ---
$lastSheet = $wb.Sheets.Item($wb.Sheets.Count)
$skate = $wb.Worksheets.Add($null, $lastSheet)
$skate.Name = "Support de skate"

$skate.Range("A1").Value = "Coût (€)"
$skate.Range("B1").Value = 12

$skate.Range("A2").Value = "Resistance (MPa)"
$skate.Range("B2").Value = 110

$skate.Range("A3").Value = "Temps de cycle moyen (min) "
$skate.Range("B3").Value = 5

$skate.Range("A4").Value = "Gamme fabrication"
$skate.Range("B4").Value = "Nom operation"
$skate.Range("C4").Value = "Machine"
$skate.Range("D4").Value = "Temps fabrication"

$skate.Range("A5").NumberFormat = "@"
$skate.Range("A5").Value = "1"
$skate.Range("B5").Value = "Moulage"
$skate.Range("C5").Value = "Moulage automatique"
$skate.Range("D5").NumberFormat = "@"
$skate.Range("D5").Value = "5"

$skate.Range("A6").NumberFormat = "@"
$skate.Range("A6").Value = "2"
$skate.Range("B6").Value = "Usinage"
$skate.Range("C6").Value = "Fraiseuse 3 axes manuelle"
$skate.Range("D6").NumberFormat = "@"
$skate.Range("D6").Value = "5"

$skate.Range("A7").NumberFormat = "@"
$skate.Range("A7").Value = "3"
$skate.Range("B7").Value = "Finition"
$skate.Range("C7").Value = "Fraiseuse 3 axes manuelle"
$skate.Range("D7").NumberFormat = "@"
$skate.Range("D7").Value = "8"

Write-Output "edit applied"
